$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data
# Using a leading apostrophe forces text storage (so single-decimal numeric-
# looking prices like "217.29" are not auto-converted to Number type), then
# resetting Style back to Normal strips the quote-prefix style Excel applies,
# keeping the cell style identical to the untouched original (no "s" attr).
$ws.Range("D2").Value = "'27.212.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.20%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.645.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.09%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'217.29"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "'  +1.91%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +1.35%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.12%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.84%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.29%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.876.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.15%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.649.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.32%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  +0.30%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  +2.69%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'67.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.11%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'27.193.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.13%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +1.11%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'219.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.50%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.06%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +5.71%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.96%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.35%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +0.54%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'147.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.40%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +2.62%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E28").Value = "'  -0.03%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.58%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.56%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.08%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +0.28%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +1.30%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +1.88%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.262.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.08%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.31%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +1.98%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.545"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.56%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.69%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -0.08%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.41%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +1.76%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.786.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.05%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'61.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.62%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'91.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.80%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +2.47%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.13%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.32%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +0.10%  "
$ws.Range("E51").Style = "Normal"
